# This script applies the IFRS financial-data correction described in the
# commit "error solve ifrs list": the per-row figures for 유안타증권 in
# rows 2-9 (columns D:AJ) are replaced with corrected values. The "FCF"
# column (U) is dropped entirely for rows 2-6, and several trailing metric
# columns are cleared where the corrected data no longer has values
# (AG/AH in row 5, AG/AH/AI in row 6). Rows 7-9 no longer carry any of the
# yearly estimate figures (columns D:AJ are cleared), leaving only the
# identifying columns A:C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")

$data = @{}
$data[2] = @(8480,-1149,-1149,-1683,-1695,-1698,3,71498,62302,9196,9176,20,10456,-4253,883,7995,28,$null,8832,-13.55,-19.99,-18.72,-2.59,677.5,-6.67,-956,-3.87,4540,0.8100000000000001,0,0,0,196214395)
$data[3] = @(12174,220,220,479,581,577,4,89232,79335,9897,9897,1,10624,-10008,960,11550,84,$null,14018,1.81,4.78,6.05,0.72,801.61,-1.54,273,13.47,4816,0.76,0,0,0,199577435)
$data[4] = @(14605,132,132,330,313,313,0,103631,93303,10327,10327,0,10624,-9905,-206,9504,49,$null,16332,0.9,2.14,3.1,0.33,903.47,2.51,147,21.27,5026,0.62,0,0,0,199577435)
$data[5] = @(19280,585,585,713,707,707,0,116512,105517,10994,10994,0,10624,-2663,198,904,27,$null,19929,3.03,3.67,6.63,0.64,959.75,8.789999999999999,333,11.6,5351,0.72,$null,$null,0,199577436)
$data[6] = @(20372,911,911,1142,1047,1047,$null,118618,106723,11895,11895,$null,10625,728,854,-2008,40,$null,18109,4.47,5.14,9.15,0.89,897.1799999999999,17.26,493,6.37,5789,0.54,$null,$null,$null,199596576)
$data[7] = @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
$data[8] = @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
$data[9] = @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null)

foreach ($r in $data.Keys) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = $cols[$i] + $r
        $value = $rowValues[$i]
        if ($null -eq $value) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
